$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "_tejgtotfun_f5amb"
$ws.Range("C2").Value = 0.005504846381997101
$ws.Range("B3").Value = "dfgdevpiagfun_f5ct05sanpc"
$ws.Range("C3").Value = 0.004731700195798244
$ws.Range("B4").Value = "_tejgct_r08gstcr"
$ws.Range("C4").Value = 0.002959397166986928
$ws.Range("B5").Value = "_tejgtotfun_f2opseg"
$ws.Range("C5").Value = 0.002675886849323625
$ws.Range("B6").Value = "bvleche_16"
$ws.Range("C6").Value = 0.002442280975516227
$ws.Range("B7").Value = "_tejgfun_f5r08ct05pgrcopc"
$ws.Range("C7").Value = 0.002404761835157466
$ws.Range("B8").Value = "_tejgrb_redr"
$ws.Range("C8").Value = 0.002314993921975503
$ws.Range("B9").Value = "_tejgtotfun_f5r08opseg"
$ws.Range("C9").Value = 0.002188365841349146
$ws.Range("B10").Value = "_tejgge_r08ct05biser"
$ws.Range("C10").Value = 0.002147029025418947
$ws.Range("B11").Value = "_tejgft_redr"
$ws.Range("C11").Value = 0.002057015610393463
$ws.Range("B12").Value = "_tejgfun_f5ct05pgrco"
$ws.Range("C12").Value = 0.00203007639476122
$ws.Range("B13").Value = "_tejgct_r09gstcr"
$ws.Range("C13").Value = 0.0019213393270223
$ws.Range("B14").Value = "_tejgtotfun_f2pgrco"
$ws.Range("C14").Value = 0.001918691447780595
$ws.Range("B15").Value = "tdvgkft_rdet"
$ws.Range("C15").Value = 0.001914663816750941
$ws.Range("B16").Value = "dfgdevpiagfun_f4ct06transpc"
$ws.Range("C16").Value = 0.001911906000148012
$ws.Range("B17").Value = "tejgtotfun_f2amb"
$ws.Range("C17").Value = 0.00187705051137771
$ws.Range("B18").Value = "_tejgfun_f5r18ct05opseg"
$ws.Range("C18").Value = 0.001809601496711791
$ws.Range("B19").Value = "dfgdevpiagfun_f5r18ct05sanpc"
$ws.Range("C19").Value = 0.001806101396983417
$ws.Range("B20").Value = "dfgdevpiagtotfun_f4transpc"
$ws.Range("C20").Value = 0.001788618579459514
$ws.Range("B21").Value = "_dfgdevpiagct_r18gstcp"
$ws.Range("C21").Value = 0.001695747656399214
$ws.Range("B22").Value = "_pimgfun_f5ct06opseg"
$ws.Range("C22").Value = 0.001612650621782239
$ws.Range("B23").Value = "_tejgfun_f2ct05pgrco"
$ws.Range("C23").Value = 0.001594955783620757
$ws.Range("B24").Value = "_dfgpimpiafun_f5ct06opsegpc"
$ws.Range("C24").Value = 0.001531251095061599
$ws.Range("B25").Value = "piagtotfun_f5r07protspc"
$ws.Range("C25").Value = 0.00149962816255714
$ws.Range("B26").Value = "_dfgpimpiafun_f5ct06opseg"
$ws.Range("C26").Value = 0.001451930471402063
$ws.Range("B27").Value = "_tejgge_r09ct05biser"
$ws.Range("C27").Value = 0.001432017498283682
$ws.Range("B28").Value = "_dfgpimpiagkft_rdet"
$ws.Range("C28").Value = 0.001429496717933424
$ws.Range("B29").Value = "_tejgfun_f5ct05amb"
$ws.Range("C29").Value = 0.001418861413727493
$ws.Range("B30").Value = "_dfgpimpiatotfun_f5opsegpc"
$ws.Range("C30").Value = 0.001407321363397562
$ws.Range("B31").Value = "_pimgfun_f5ct06opsegpc"
$ws.Range("C31").Value = 0.001398633022068218
$ws.Range("B32").Value = "tejgfun_f2ct05amb"
$ws.Range("C32").Value = 0.001397021039243188
$ws.Range("B33").Value = "_tejgtotfun_f5pgrco"
$ws.Range("C33").Value = 0.001394035484187648
$ws.Range("B34").Value = "_tejgfun_f1ct06pgrco"
$ws.Range("C34").Value = 0.001393500231854979
$ws.Range("B35").Value = "_pimgtotfun_f5opsegpc"
$ws.Range("C35").Value = 0.00139346022474405
$ws.Range("B36").Value = "_tejgfun_f5ct05opseg"
$ws.Range("C36").Value = 0.001343742868896097
$ws.Range("B37").Value = "_tejgge_r09ct05biserpc"
$ws.Range("C37").Value = 0.001329173390418222
$ws.Range("B38").Value = "devppimfun_f2ct06agro"
$ws.Range("C38").Value = 0.001310439158062251
$ws.Range("B39").Value = "_tejgfun_f2ct05amb"
$ws.Range("C39").Value = 0.001300114910866389
$ws.Range("B40").Value = "_tejgge_r08ct05pobso"
$ws.Range("C40").Value = 0.001297698782067145
$ws.Range("B41").Value = "pimgfun_f1ct06prots"
$ws.Range("C41").Value = 0.001290050175474327
$ws.Range("B42").Value = "_dfgdevpiagkftr18_rdet"
$ws.Range("C42").Value = 0.001270474428132654
$ws.Range("B43").Value = "_tejgfun_f5r18ct05amb"
$ws.Range("C43").Value = 0.001231091028748842
$ws.Range("B44").Value = "dfgpimpiafun_f5r07ct06protspc"
$ws.Range("C44").Value = 0.001230924628802432
$ws.Range("B45").Value = "_tejgct_r09gstcrpc"
$ws.Range("C45").Value = 0.00121609936968837
$ws.Range("B46").Value = "devppimgkft_rdet"
$ws.Range("C46").Value = 0.001197019500033859
$ws.Range("B47").Value = "_devppimct_r00gstcr"
$ws.Range("C47").Value = 0.001160932097414928
$ws.Range("B48").Value = "tejgfun_f5ct05pgrco"
$ws.Range("C48").Value = 0.001158218487768762
$ws.Range("B49").Value = "_tejgfun_f2ct05opseg"
$ws.Range("C49").Value = 0.001141764598998171
$ws.Range("B50").Value = "devppimfun_f5r18ct05cydep"
$ws.Range("C50").Value = 0.001107854731024029
$ws.Range("B51").Value = "dfgpimpiagkft_rdetpc"
$ws.Range("C51").Value = 0.001099558488742017
